# Reorders the player roster rows (A2:C19) on the active sheet and swaps
# "Jaxson Hayes / PF,C / Los Angeles Lakers" for
# "RJ Barrett / SG,SF,PF / Toronto Raptors".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Jeremy Sochan", "SF,PF", "San Antonio Spurs"),
    @("Amen Thompson", "SG,SF,PF", "Houston Rockets"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Goga Bitadze", "C", "Orlando Magic"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
